$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Daten")

function Set-DataStyle($cell) {
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.HorizontalAlignment = -4152  # xlRight
}

# --- New row 5: Hospital (written first so new shared strings "Hospital"/"True"
#     get the lower indices, matching the target string table order) ---
$ws.Cells.Item(5, 1).Value = "Hospital"
Set-DataStyle $ws.Cells.Item(5, 1)

$ws.Cells.Item(5, 2).Value = 3000
Set-DataStyle $ws.Cells.Item(5, 2)

$ws.Cells.Item(5, 3).Value = 0.005
Set-DataStyle $ws.Cells.Item(5, 3)

$ws.Cells.Item(5, 4).Value = 250
Set-DataStyle $ws.Cells.Item(5, 4)

$ws.Cells.Item(5, 5).Value = 3000
Set-DataStyle $ws.Cells.Item(5, 5)

$ws.Cells.Item(5, 6).Value = 20
Set-DataStyle $ws.Cells.Item(5, 6)

# Write "True" as literal text (not boolean) via a formula + paste-values trick
$gCell = $ws.Cells.Item(5, 7)
$gCell.Formula = '="True"'
$gCell.Copy()
$gCell.PasteSpecial(-4163)  # xlPasteValues
Set-DataStyle $gCell

$ws.Cells.Item(5, 8).Value = 0.1
Set-DataStyle $ws.Cells.Item(5, 8)

# --- New column H header ("beta") ---
$ws.Cells.Item(1, 8).Value = "beta"
Set-DataStyle $ws.Cells.Item(1, 8)

# --- New beta values for existing rows 2-4 ---
$ws.Cells.Item(2, 8).Value = 0.1
Set-DataStyle $ws.Cells.Item(2, 8)

$ws.Cells.Item(3, 8).Value = 0.05
Set-DataStyle $ws.Cells.Item(3, 8)

$ws.Cells.Item(4, 8).Value = 0.1
Set-DataStyle $ws.Cells.Item(4, 8)

# --- Update the saved selection to H5 ---
$ws.Range("H5").Select()
